# Contenido hasta semana 4
# Fill in the GitHub usernames ("Usuario de GitHub", column F) for the
# students that did not have one registered yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = "Camilo-AndradePerez"
$ws.Range("F7").Value  = "LauraCaballero14"
$ws.Range("F13").Value = "maflecha"
$ws.Range("F15").Value = "javierfranco27"
$ws.Range("F29").Value = "Arpenahi"
$ws.Range("F33").Value = "lipinilla685"
$ws.Range("F34").Value = "HolguerRangel"
$ws.Range("F36").Value = "mrincon19"
$ws.Range("F38").Value = "stephanierojas1234"
$ws.Range("F40").Value = "LauraSanchez9585"

$ws.Range("E2").Select()
